# Apply "Optuna Attempt (go back with original)" data updates.
$wb = $excel.ActiveWorkbook

# --- Sheet: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("D2").Value = 31
$ws1.Range("H2").Value = 17.77
$ws1.Range("L2").Value = 0.9399999999999999

$ws1.Range("D3").Value = 30
$ws1.Range("H3").Value = 17.33

$ws1.Range("H4").Value = 18.85
$ws1.Range("L4").Value = 1.15

$ws1.Range("D5").Value = 27
$ws1.Range("H5").Value = 17.19
$ws1.Range("L5").Value = 0.95

$ws1.Range("H6").Value = 16.81
$ws1.Range("L6").Value = 1.15

$ws1.Range("H7").Value = 16.64
$ws1.Range("L7").Value = 1.02

$ws1.Range("H8").Value = 14.15
$ws1.Range("L8").Value = 1.14

$ws1.Range("D9").Value = 27
$ws1.Range("H9").Value = 13.3
$ws1.Range("L9").Value = 1.15

$ws1.Range("H10").Value = 14.19
$ws1.Range("L10").Value = 1.09

$ws1.Range("D11").Value = 25
$ws1.Range("H11").Value = 12.34
$ws1.Range("L11").Value = 1.09

$ws1.Range("H12").Value = 10.91
$ws1.Range("L12").Value = 0.93

$ws1.Range("H13").Value = 10.43
$ws1.Range("L13").Value = 0.95

$ws1.Range("D14").Value = 25
$ws1.Range("H14").Value = 9.32
$ws1.Range("L14").Value = 0.9399999999999999

$ws1.Range("H15").Value = 8.66
$ws1.Range("L15").Value = 1.14

$ws1.Range("H16").Value = 7.86
$ws1.Range("L16").Value = 1.02

$ws1.Range("H17").Value = 6.86
$ws1.Range("L17").Value = 0.8100000000000001

# --- Sheet: "Summary" ---
# Column B on this sheet stores numbers as text (inline strings), so
# prefix with an apostrophe to force text entry and keep the "General"
# number format intact (matches original file layout).
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'414"
$ws2.Range("B10").Value = "'219"
$ws2.Range("B11").Value = "'114"
$ws2.Range("B12").Value = "'31"
